$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

$ws.Range("E8:E21").Value = "no"

$ws.Range("E8").Select()
